# "mean and median push" - update simulated cluster-oversight values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{ Row=2;  A="AON"; B="MMC"; C=0.9796450739087524;  D=0.120720697581247;   E=0.9048669289296306; F=0.9922189188275892; G=0.4307694124029075;  H=0.06260296235450927 },
    @{ Row=3;  A="AON"; B="BRO"; C=0.9623298439021065;  D=0.1028338431253394;  E=0.9048669289296306; F=0.9975325979692565; G=0.7975044194634193;  H=0.3326450380065146 },
    @{ Row=4;  A="AJG"; B="MMC"; C=0.9783479380831384;  D=0.2600292245963809;  E=0.9948012752887331; F=0.9922189188275892; G=0.1903475230646941;  H=0.2186635195045711 },
    @{ Row=5;  A="BRO"; B="MMC"; C=0.9777233743065304;  D=0.1179100127470729;  E=0.9975325979692565; F=0.9922189188275892; G=0.9587888911066282;  H=0.141549606418864  },
    @{ Row=6;  A="AJG"; B="AON"; C=0.9657091117319888;  D=0.100073378760585;   E=0.9948012752887331; F=0.9048669289296306; G=0.3714655397054613;  H=0.3852252068075139 },
    @{ Row=7;  A="AJG"; B="BRO"; C=0.9771805473804041;  D=0.2982061204769986;  E=0.9948012752887331; F=0.9975325979692565; G=0.9681175598233733;  H=0.224775862899768  },
    @{ Row=8;  A="AJG"; B="WTW"; C=0.9241385728782294;  D=0.7330819931126087;  E=0.9948012752887331; F=0.6913029326157746; G=0.165423191350179;   H=0.8990479535705498 },
    @{ Row=9;  A="MMC"; B="WTW"; C=0.8981858138266483;  D=0.8417653326866769;  E=0.9922189188275892; F=0.6913029326157746; G=0.1817159148413596;  H=0.8661209308406749 },
    @{ Row=10; A="BRO"; B="WTW"; C=0.9097472141186467;  D=0.5355755997733355;  E=0.9975325979692565; F=0.6913029326157746; G=0.3779119866090953;  H=0.9347848457383896 },
    @{ Row=11; A="AON"; B="WTW"; C=0.9189908896166797;  D=0.6616385184807054;  E=0.9048669289296306; F=0.6913029326157746; G=0.8492948757379242;  H=0.667027911391471  }
)

foreach ($row in $data) {
    $r = $row.Row
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Cells.Item($r, 7).Value = $row.G
    $ws.Cells.Item($r, 8).Value = $row.H
}
